# Fehler in Test-Excel verbaut
#
# The "Jugend" sheet had a stray header row accidentally left in place
# (row 1 duplicated the table header that every other sheet already
# carries via its own header row / styling). Remove it, which shifts
# all the real data rows up by one. Also make "Jugend" the active /
# selected sheet (it was "GS33-35" before).

$wb = $excel.ActiveWorkbook

# Jump to the sheet that has the stray row and make it the active tab.
$ws = $wb.Worksheets.Item("Jugend")
$ws.Activate()

# Delete the duplicate header row; everything below shifts up by one
# and formulas/relative refs are renumbered automatically.
$ws.Rows(1).Delete()

# Leave the selection where it ended up after the edit.
$ws.Range("A7").Select()
